# Update odds values on Sheet1 (rows 3-21) to match the 2024-12-05 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 2.63
$ws.Range("X3").Value = 17
$ws.Range("AD3").Value = 9
$ws.Range("AF3").Value = 34
$ws.Range("AG3").Value = 101
$ws.Range("AI3").Value = 21
$ws.Range("AM3").Value = 23
$ws.Range("G4").Value = 3.6
$ws.Range("I4").Value = 2.05
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("W4").Value = 11
$ws.Range("X4").Value = 19
$ws.Range("AI4").Value = 9.5
$ws.Range("AK4").Value = 17
$ws.Range("AO4").Value = 21
$ws.Range("AT4").Value = 2.75
$ws.Range("AY4").Value = 11
$ws.Range("AZ4").Value = 21
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 2.63
$ws.Range("X5").Value = 9
$ws.Range("AH5").Value = 11
$ws.Range("AL5").Value = 34
$ws.Range("AP5").Value = 21
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
$ws.Range("H9").Value = 3.5
$ws.Range("K9").Value = 2.05
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.6
$ws.Range("AC9").Value = 7.5
$ws.Range("AF9").Value = 67
$ws.Range("AP9").Value = 23
$ws.Range("AS9").Value = 201
$ws.Range("AV9").Value = 67
$ws.Range("BA9").Value = 101
$ws.Range("G10").Value = 1.95
$ws.Range("H10").Value = 2.9
$ws.Range("I10").Value = 4.75
$ws.Range("J10").Value = 2.75
$ws.Range("K10").Value = 1.83
$ws.Range("L10").Value = 5.5
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("O10").Value = 1.62
$ws.Range("P10").Value = 2.2
$ws.Range("S10").Value = 1.67
$ws.Range("T10").Value = 2.1
$ws.Range("X10").Value = 7.5
$ws.Range("Z10").Value = 17
$ws.Range("AA10").Value = 21
$ws.Range("AD10").Value = 6
$ws.Range("AH10").Value = 9
$ws.Range("AI10").Value = 21
$ws.Range("AJ10").Value = 19
$ws.Range("AK10").Value = 51
$ws.Range("AN10").Value = 3.6
$ws.Range("AO10").Value = 12
$ws.Range("AT10").Value = 2.1
$ws.Range("AU10").Value = 10
$ws.Range("AX10").Value = 6
$ws.Range("AY10").Value = 29
$ws.Range("BA10").Value = 126
$ws.Range("G21").Value = 8.5
$ws.Range("I21").Value = 1.33
$ws.Range("J21").Value = 8
$ws.Range("K21").Value = 2.5
$ws.Range("Q21").Value = 1.67
$ws.Range("R21").Value = 2.15
$ws.Range("S21").Value = 1.3
$ws.Range("T21").Value = 3.4
$ws.Range("U21").Value = 2
$ws.Range("V21").Value = 1.75
$ws.Range("AA21").Value = 51
$ws.Range("AC21").Value = 13
$ws.Range("AE21").Value = 21
$ws.Range("AG21").Value = 351
$ws.Range("AH21").Value = 7
$ws.Range("AI21").Value = 6.5
$ws.Range("AJ21").Value = 8.5
$ws.Range("AK21").Value = 8.5
$ws.Range("AL21").Value = 11
$ws.Range("AM21").Value = 26
$ws.Range("AQ21").Value = 151
$ws.Range("AR21").Value = 151
$ws.Range("AS21").Value = 301
$ws.Range("AT21").Value = 3.4
$ws.Range("AV21").Value = 51
$ws.Range("AX21").Value = 3.25
$ws.Range("AY21").Value = 6.5
$ws.Range("AZ21").Value = 17
$ws.Range("BC21").Value = 126
